# Update "想去人数" (want-to-go count) figures in column F, reflecting a
# refreshed scrape of the 苏州-漫展信息 event listing (gh-pages output at
# commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 37
$ws1.Range("F3").Value  = 105
$ws1.Range("F4").Value  = 1510
$ws1.Range("F5").Value  = 212
$ws1.Range("F6").Value  = 46
$ws1.Range("F7").Value  = 350
$ws1.Range("F8").Value  = 9914
$ws1.Range("F10").Value = 122
$ws1.Range("F14").Value = 6854
$ws1.Range("F16").Value = 636
$ws1.Range("F18").Value = 201

# --- Sheet "全部类型" (all types, superset incl. 演出 rows) ------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 37
$ws4.Range("F3").Value  = 105
$ws4.Range("F4").Value  = 1510
$ws4.Range("F5").Value  = 212
$ws4.Range("F7").Value  = 46
$ws4.Range("F8").Value  = 350
$ws4.Range("F11").Value = 9914
$ws4.Range("F13").Value = 122
$ws4.Range("F17").Value = 6855
$ws4.Range("F19").Value = 636
$ws4.Range("F21").Value = 201
